$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new soil texture, "loam", was added to the lookup list in column B,
# just below "silt"/"silt loam". Shift the existing entries (sand loam,
# loamy sand, fine sand, med sand, coarse sand, fine gravel, med gravel,
# coarse gravel, cobble) down by one row to make room, writing the values
# directly so the row formatting (column A/B styles) and column C index
# values are left untouched.
$textures = @("sand loam", "loamy sand", "fine sand", "med sand", "coarse sand", "fine gravel", "med gravel", "coarse gravel", "cobble")
for ($i = $textures.Length - 1; $i -ge 0; $i--) {
    $ws.Cells.Item(5 + $i, 2).Value = $textures[$i]
}
$ws.Cells.Item(4, 2).Value = "loam"

# "cobble" now lands on row 13, which previously had no "texture grade"
# value; continue the existing 1,2,3... sequence from column C.
$ws.Cells.Item(13, 3).Value = 12

# Extend the table with one more blank row at the bottom (row 26), matching
# the formatting already used by row 25.
$ws.Range("A25:B25").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(26, 1).ClearContents()
$ws.Cells.Item(26, 2).ClearContents()
$excel.CutCopyMode = $false

# Record the cell that was selected at the end of the editing session.
$ws.Range("E10").Select()
